# Generate Report for Handback
# The handback file for this resource was (re)generated, so the status
# flips from "In Translation" to "Handed back: in sync with en-US", the
# "Latest Handback DateTime" timestamps for each locale move forward to
# the new handback time, and the stale "handback not latest" error is
# cleared now that the handback report is fresh.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item(1)   # "Overview"
$zhcn     = $wb.Worksheets.Item(2)   # "zh-cn"
$dede     = $wb.Worksheets.Item(3)   # "de-de"

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: Status column per locale (E = zh-cn, F = de-de) ---
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = 29.2
$overview.Columns.Item(6).ColumnWidth = 29.2

# --- zh-cn detail sheet ---
$zhcn.Range("C2").Value = $statusText            # Status
$zhcn.Range("K2").Value = "2016-11-15 16:19:01"  # Latest Handback DateTime
$zhcn.Range("P2").Value = ""                     # Error Detail cleared
$zhcn.Columns.Item(3).ColumnWidth = 29.2
$zhcn.Columns.Item(16).ColumnWidth = 12.9

# --- de-de detail sheet ---
$dede.Range("C2").Value = $statusText            # Status
$dede.Range("K2").Value = "2016-11-15 16:19:19"  # Latest Handback DateTime
$dede.Range("P2").Value = ""                     # Error Detail cleared
$dede.Columns.Item(3).ColumnWidth = 29.2
$dede.Columns.Item(16).ColumnWidth = 12.9
